$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.48
$ws.Range("I2").Value = 6
$ws.Range("U2").Value = 1.62
$ws.Range("V2").Value = 2.2
$ws.Range("Y2").Value = 9.5
$ws.Range("Z2").Value = 13
$ws.Range("AF2").Value = 41
$ws.Range("AG2").Value = 151
$ws.Range("AK2").Value = 67
$ws.Range("AU2").Value = 7.5
$ws.Range("AX2").Value = 7.5
$ws.Range("BA2").Value = 81
$ws.Range("BB2").Value = 81
$ws.Range("G3").Value = 2.1
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 3.7
$ws.Range("J3").Value = 2.88
$ws.Range("L3").Value = 4.33
$ws.Range("M3").Value = 1.08
$ws.Range("O3").Value = 1.4
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 1.95
$ws.Range("V3").Value = 1.8
$ws.Range("X3").Value = 9.5
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 19
$ws.Range("AC3").Value = 8
$ws.Range("AG3").Value = 351
$ws.Range("AH3").Value = 9.5
$ws.Range("AI3").Value = 17
$ws.Range("AJ3").Value = 13
$ws.Range("AK3").Value = 41
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 12
$ws.Range("AT3").Value = 2.5
$ws.Range("AV3").Value = 67
$ws.Range("AX3").Value = 5.5
$ws.Range("AY3").Value = 21
$ws.Range("BB3").Value = 101
$ws.Range("BD3").Value = 151
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
